$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default, unstyled body cell) used to avoid leaving a
# stray "Text" number-format style behind on cells we must force to text
# so that numeric-looking strings (e.g. "1.250") keep their exact text value.
$refStyle = $ws.Range("C2").Style

$ws.Range("D2").Value = "30.148.79"
$ws.Range("E2").Value = "  +5.72%  "
$ws.Range("D3").Value = "1.920.60"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.09"
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = "  +4.64%  "
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5219"
$ws.Range("D7").Style = $refStyle
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4089"
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = "  +4.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08547"
$ws.Range("D9").Style = $refStyle
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.128"
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.37"
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = "  +9.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.419"
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = "  +3.22%  "
$ws.Range("D14").Value = "1.913.76"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.412"
$ws.Range("D15").Style = $refStyle
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.43"
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = "  +4.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001115"
$ws.Range("D18").Style = $refStyle
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.42"
$ws.Range("D20").Style = $refStyle
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("D21").Style = $refStyle
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.015"
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("D23").Value = "30.149.16"
$ws.Range("E23").Value = "  +5.62%  "
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.211"
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "2.136.36"
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.10"
$ws.Range("D27").Style = $refStyle
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.84"
$ws.Range("D28").Style = $refStyle
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.448"
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.24"
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.083"
$ws.Range("D31").Style = $refStyle
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1063"
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.052"
$ws.Range("D33").Style = $refStyle
$ws.Range("E33").Value = "  +5.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.635"
$ws.Range("D34").Style = $refStyle
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02496"
$ws.Range("D35").Style = $refStyle
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06611"
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2208"
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.235"
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = "  +4.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.182"
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = "  +2.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.925"
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6554"
$ws.Range("D41").Style = $refStyle
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.253"
$ws.Range("D42").Style = $refStyle
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("E43").Value = "  +4.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6161"
$ws.Range("D44").Style = $refStyle
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.25"
$ws.Range("D45").Style = $refStyle
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.763"
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.081"
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = "  +3.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.250"
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.63"
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.169"
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = "  +10.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.83"
$ws.Range("D51").Style = $refStyle
$ws.Range("E51").Value = "  +4.25%  "
